$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# The document opens with two short title paragraphs:
#   1) "Kelly Cash"
#   2) "Excel Challenge Homework Report"
# They need to become a single paragraph reading
#   "Excel Kickstarter Campaign Report"
# split across three runs: "Excel " | "Kickstarter Campaign" | " Report"
# ---------------------------------------------------------------------

# 1) Turn "Kelly Cash" into "Excel " (keeps paragraph 1's own run).
$d.Content.Find.Execute("Kelly Cash", $true, $false, $false, $false, $false, $true, 1, $false, "Excel ", 2) | Out-Null

# 2) Turn "Excel Challenge Homework Report" into "Kickstarter Campaign Report".
$d.Content.Find.Execute("Excel Challenge Homework", $true, $false, $false, $false, $false, $true, 1, $false, "Kickstarter Campaign", 2) | Out-Null

# 3) Split paragraph 2 ("Kickstarter Campaign Report") right before " Report"
#    so the two pieces end up as independent runs once merged below.
$p2 = $d.Paragraphs(2).Range
$splitRange = $d.Range($p2.Start, $p2.End)
$splitRange.Find.Execute(" Report") | Out-Null
$breakPoint = $d.Range($splitRange.Start, $splitRange.Start)
$breakPoint.InsertParagraphAfter()

# 4) Move the ("Kickstarter Campaign") paragraph's text into paragraph 1,
#    right before its end-of-paragraph mark, as its own run, then drop
#    the now-empty paragraph.
$piece = $d.Paragraphs(2).Range
$pieceText = $d.Range($piece.Start, $piece.End - 1)
$pieceText.Cut()
$dest = $d.Range($d.Paragraphs(1).Range.End - 1, $d.Paragraphs(1).Range.End - 1)
$dest.Paste()
$d.Paragraphs(2).Range.Delete()

# 5) Do the same for the (" Report") paragraph.
$piece2 = $d.Paragraphs(2).Range
$piece2Text = $d.Range($piece2.Start, $piece2.End - 1)
$piece2Text.Cut()
$dest2 = $d.Range($d.Paragraphs(1).Range.End - 1, $d.Paragraphs(1).Range.End - 1)
$dest2.Paste()
$d.Paragraphs(2).Range.Delete()
